$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.934.06"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "2.772.37"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.66%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.113"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.394"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").Value = "3.270.57"
$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15").Value = "63.868.13"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000157"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.78%  "

$ws.Range("D17").Value = "2.786.28"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.574"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.175"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").Value = "0.0₃0947"
$ws.Range("E27").Value = "  +9.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "170.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.60%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "336.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0605"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.646"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0260"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.16%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.148.91"
$ws.Range("E51").Value = "  -0.59%  "
